$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 501; $r++) {
    $ws.Cells.Item($r, 9).Formula = "=(`$L`$2/COUNT(`$B`$2:`$B`$501)) * (H$r / A$r)"
}

$ws.Range("I2:I501").Select()
